$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.243.62"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.17%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.866.55"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.08%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7184"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.21%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "240.80"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +1.27%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.22%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07825"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.04%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3083"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.67%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "25.12"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.47%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08263"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.37%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.884.56"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.14%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.7210"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.18%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.224"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.03%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "90.60"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.69%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "29.276.15"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.14%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.858"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.91%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "243.63"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.70%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007798"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.02%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.125.53"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.83%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.16"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.36%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.002"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.24%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.981"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +5.44%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.001"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1602"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +10.45%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "162.28"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.04%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.934"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.19%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.21"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.45%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.341"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.99%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.494"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.41%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.397"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.95%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.097"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.18%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05203"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.36%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.925"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.76%  "

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.07%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7279"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.45%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.683"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.36%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01851"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.11%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.693"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.07%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.168.47"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.00%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9025"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.47%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.115"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.86%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "72.49"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.96%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.002"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.23%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "101.83"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.16%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.022.25"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.66%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5286"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.07%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.776"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.34%  "

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +4.04%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.302"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.24%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.865"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +5.03%  "
